$wb = $excel.ActiveWorkbook

# Add new worksheet at the end, after "20201028"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "20201031"

# Fill header row
$newSheet.Range("A1").Value = "ID"
$newSheet.Range("B1").Value = "Shop ID"
$newSheet.Range("C1").Value = "SQL"
$newSheet.Range("A1:C1").Font.Size = 12
$newSheet.Range("A1:C1").Font.Color = 0

# Fill data rows 2-7 (ID 1-6), shop id constant
$shopId = "32fe0cfd-0254-11eb-ba65-065a10bcba76"
for ($i = 1; $i -le 6; $i++) {
    $row = $i + 1
    $newSheet.Cells.Item($row, 1).Value = $i
    $newSheet.Cells.Item($row, 2).Value = $shopId
    $newSheet.Cells.Item($row, 3).Formula = "=_xlfn.CONCAT(""INSERT INTO photos(restaurant_id, name, type) VALUES(UuidToBin('"", B$row, ""'), LPAD("", A$row, "", 7, '0'), 'dish'"", "");"")"
}

$wb.Worksheets.Item("20201028").Range("C2").Select()

$newSheet.Activate()
$newSheet.Range("F11").Select()
